$d = $word.ActiveDocument

# 1) Update the initiation/report date field from 06/23/2015 to 07/07/2015
#    (appears twice in the document, in table cells labeled
#    "2. INITIATION DATE:" and "7. DATE:")
$d.Content.Find.Execute("06/23/2015", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "07/07/2015", 2)

# 2) Update the verification "Estimated Completion Date" from June 4 2015
#    to July 1 2015 (appears three times, once per verification block)
$d.Content.Find.Execute("June 4 2015", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "July 1 2015", 2)
